$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update timing values in rows 5-7
$ws.Range("B5").Value = 0.001062154769897461
$ws.Range("B6").Value = 0.0006451606750488281
$ws.Range("B7").Value = 0.008590936660766602

# 2. Convert tuple-style text to list-style text (parentheses -> brackets)
$ws.Range("A8").Value = "[[3, 3], [3, 2], [3, 0], [2, 0], [2, 2], [2, 3], [0, 3], [0, 2], [0, 0], [2, 1], [1, 2], [1, 0], [1, 1], [0, 1]]"
$ws.Range("A49").Value = "[[3, 3], [3, 2], [3, 0], [2, 0], [2, 1], [1, 1], [1, 2], [0, 2], [0, 3], [1, 3], [0, 1], [2, 2], [1, 0], [3, 1]]"
$ws.Range("A111").Value = "[[3, 3], [2, 3], [1, 2], [2, 0], [1, 0], [0, 2], [0, 1], [0, 0], [1, 1], [2, 1], [2, 2], [1, 3], [3, 2], [3, 1]]"
$ws.Range("A169").Value = "[[0, 2], [1, 2], [1, 1], [2, 1], [1, 0], [0, 1], [2, 2], [1, 3], [0, 3], [0, 0], [2, 0], [2, 3], [3, 2], [3, 1]]"
$ws.Range("A219").Value = "[[0, 2], [0, 1], [1, 1], [0, 0], [0, 3], [2, 3], [2, 0], [3, 0], [3, 2], [1, 2], [1, 0], [1, 3], [2, 1], [2, 2]]"

# 3. Insert a new row at 253 (pushes old rows 253-257 down to 254-258)
$ws.Rows.Item(253).Insert()
$ws.Range("A253").Value = "move_fidelity"
$ws.Range("B253").Value = 0.9978635150727226

# 4. Update the new "total time:" row (now row 257) with the new value
$ws.Range("B257").Value = 0.03424215316772461

$wb.Save()
Write-Output "done"
